# Applies:
#  1) Straight Arrow Connector 34 (slide 2, the "Trekke tilbake avtale" arrow):
#     tailEnd none -> arrow (so it points both ways).
#  2) TextBox 69 (slide 2, the "avslutte avtale" bullet box): shrink height,
#     drop the "Hvis avtalen avsluttes ..." bullet, and change "Ny
#     avtalestatus overføres til Digipost" to "Ny avtalestatus overføres til
#     den andre parten".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate shapes by their (unique) original position rather than a bare
# positional index, since a couple of shapes on this slide share a Name.
$arrow = $null
$box = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Straight Arrow Connector 34" -and [Math]::Round($sh.Left, 2) -eq 730.63) {
        $arrow = $sh
    }
    if ($sh.Name -eq "TextBox 69" -and [Math]::Round($sh.Left, 2) -eq 675.6) {
        $box = $sh
    }
}

# ---------------------------------------------------------------------------
# 1) Make the cancellation-status arrow double-headed (msoArrowheadOpen = 3,
#    matching the existing headEnd="arrow").
# ---------------------------------------------------------------------------
$arrow.Line.EndArrowheadStyle = 3

# ---------------------------------------------------------------------------
# 2) Edit the textbox that lists the "avslutte avtale" bullets.
# ---------------------------------------------------------------------------
$tr = $box.TextFrame.TextRange

$full = $tr.Text
$cr1 = $full.IndexOf([char]13)
$cr2 = $full.IndexOf([char]13, $cr1 + 1)

# Remove the whole middle paragraph ("Ny avtalestatus overføres til Digipost"),
# keeping the trailing paragraph (which carries the textbox's real
# end-of-paragraph mark) so PowerPoint doesn't leave a stray empty paragraph
# behind when the deletion reaches the end of the text body.
$midStart = $cr1 + 2
$midLen = $cr2 - $cr1
$tr.Characters($midStart, $midLen).Delete()

# The former 3rd paragraph ("Hvis avtalen avsluttes ... (6)") is now the 2nd
# paragraph; replace its whole content with the new bullet text.
$full2 = $tr.Text
$cr1b = $full2.IndexOf([char]13)
$secondStart = $cr1b + 2
$secondLen = $full2.Length - $secondStart + 1
$tr.Characters($secondStart, $secondLen).Text = "Ny avtalestatus overføres til den andre parten"

# Split that paragraph's single run into two runs ("...overføres til " / "den
# andre parten") by re-asserting the (unchanged) font size on the tail, which
# forces a run split without altering any formatting.
$full3 = $tr.Text
$splitStart = $full3.IndexOf("den andre parten") + 1
$splitLen = $full3.Length - $splitStart + 1
$tail = $tr.Characters($splitStart, $splitLen)
$tail.Font.Size = $tail.Font.Size

# Shrink the textbox now that it holds one bullet fewer (1526886 x 1100301 EMU).
$box.Height = 1100301 / 12700 + 0.00003
